# The commit swaps the content of ppt/theme/theme1.xml ("Office Theme")
# and ppt/theme/theme2.xml ("Integral") so that the design actually
# applied to the slides (the deck's single SlideMaster points at
# theme2.xml) changes from the "Integral" palette to the "Office Theme"
# palette - i.e. the presentation's Design/Theme color scheme is switched
# from Integral to the (built-in) Office colour scheme.
#
# Through the PowerPoint object model this is expressed as editing the
# twelve slots of the active design's ThemeColorScheme (dk1/lt1/dk2/lt2/
# accent1-6/hlink/folHlink) on the (single) SlideMaster, which is what a
# user does when switching the colour variant from the Design tab.

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

# Target palette = the Office Theme colours that used to live in
# theme1.xml, expressed as COM "RGB" longs (0xBBGGRR, i.e. blue/green/red
# byte order) because PowerPoint's ColorFormat.RGB is stored that way.
$cs.Item(1).RGB  = 0x000000   # dk1      = 000000
$cs.Item(2).RGB  = 0xFFFFFF   # lt1      = FFFFFF
$cs.Item(3).RGB  = 0x6A5444   # dk2      = 44546A
$cs.Item(4).RGB  = 0xE6E6E7   # lt2      = E7E6E6
$cs.Item(5).RGB  = 0xD59B5B   # accent1  = 5B9BD5
$cs.Item(6).RGB  = 0x317DED   # accent2  = ED7D31
$cs.Item(7).RGB  = 0xA5A5A5   # accent3  = A5A5A5
$cs.Item(8).RGB  = 0x00C0FF   # accent4  = FFC000
$cs.Item(9).RGB  = 0xC47244   # accent5  = 4472C4
$cs.Item(10).RGB = 0x47AD70   # accent6  = 70AD47
$cs.Item(11).RGB = 0xC16305   # hlink    = 0563C1
$cs.Item(12).RGB = 0x724F95   # folHlink = 954F72

# Best-effort: try to also restore the scheme / theme display names
# ("Office" / "Office Theme") on hosts that support renaming; harmless
# no-op where the host treats Name as read-only. (Deliberately NOT
# touching Presentation.Designs.Item(1).Name - on this host that maps to
# the slide master's <p:cSld name> attribute, which must stay untouched.)
try { $cs.Name = "Office" } catch {}
try { $p.SlideMaster.Theme.Name = "Office Theme" } catch {}
